$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.364.32"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "2.066.28"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("D5").Value = "'234.72"
$ws.Range("E5").Value = "  -1.09%  "

$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'57.28"
$ws.Range("E8").Value = "  -0.91%  "

$ws.Range("D9").Value = "'0.389"
$ws.Range("E9").Value = "  +2.04%  "

$ws.Range("D10").Value = "'0.0766"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("D11").Value = "'0.103"
$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("D12").Value = "2.368.96"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").Value = "'14.47"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").Value = "'20.59"
$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("D15").Value = "'0.775"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "'5.16"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "2.066.30"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").Value = "37.294.88"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "'6.41"
$ws.Range("E19").Value = "  +3.76%  "

$ws.Range("D20").Value = "'69.59"
$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("D21").Value = "0.0₃0814"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("D22").Value = "'227.14"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -2.45%  "

$ws.Range("D26").Value = "'166.47"
$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'1.46"
$ws.Range("E27").Value = "  +3.61%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'8.79"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").Value = "'19.11"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  +0.41%  "

$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("D32").Value = "'4.51"
$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").Value = "'0.0620"
$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("E34").Value = "  +3.01%  "

$ws.Range("E35").Value = "  -4.40%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.79"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Value = "'3.25"
$ws.Range("E38").Value = "  -2.47%  "

$ws.Range("D39").Value = "'5.50"
$ws.Range("E39").Value = "  -5.57%  "

$ws.Range("E40").Value = "  -0.89%  "

$ws.Range("D41").Value = "1.472.22"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "'4.38"
$ws.Range("E42").Value = "  -7.45%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0944"
$ws.Range("E43").Value = "  -2.59%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'96.51"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0214"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.17"
$ws.Range("E46").Value = "  +3.08%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'15.29"
$ws.Range("E47").Value = "  -4.28%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.02"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "'7.20"
$ws.Range("E49").Value = "  -0.82%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'44.75"
$ws.Range("E51").Value = "  +2.00%  "
